$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string header text updates (Volume number, week-of dates) ---
$ws.Range("A8").Value = "Volume 31   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/8/2024  Through  1/14/2024"

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 75
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 7
$ws.Range("K16").Value = 71.428571428571
$ws.Range("L16").Value = 9.090909090909
$ws.Range("M16").Value = 20
$ws.Range("N16").Value = -60
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = 16.666666666666
$ws.Range("I17").Value = 10
$ws.Range("J17").Value = 14
$ws.Range("K17").Value = -28.571428571428
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = -52.380952380952
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 4
$ws.Range("I18").Value = 2
$ws.Range("L18").Value = -66.666666666666
$ws.Range("N18").Value = -92.857142857142
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -22.222222222222
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 22.222222222222
$ws.Range("I19").Value = 14
$ws.Range("J19").Value = 15
$ws.Range("K19").Value = -6.666666666666
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 75
$ws.Range("N19").Value = 27.272727272727
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 175
$ws.Range("I20").Value = 4
$ws.Range("L20").Value = 33.333333333333
$ws.Range("M20").Value = 300
$ws.Range("N20").Value = -69.230769230769
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 5.882352941176
$ws.Range("F21").Value = 100
$ws.Range("G21").Value = 73
$ws.Range("H21").Value = 36.986301369863
$ws.Range("I21").Value = 44
$ws.Range("J21").Value = 37
$ws.Range("K21").Value = 18.918918918918
$ws.Range("L21").Value = 2.325581395348
$ws.Range("M21").Value = 62.962962962963
$ws.Range("N21").Value = -57.281553398058
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 4
$ws.Range("M22").Value = 100
$ws.Range("C23").Value = 3
$ws.Range("F23").Value = 14
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 55.555555555555
$ws.Range("I23").Value = 9
$ws.Range("K23").Value = 125
$ws.Range("M23").Value = 200
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -5.555555555555
$ws.Range("F24").Value = 73
$ws.Range("G24").Value = 67
$ws.Range("H24").Value = 8.955223880597
$ws.Range("I24").Value = 32
$ws.Range("J24").Value = 31
$ws.Range("K24").Value = 3.225806451612
$ws.Range("L24").Value = 60
$ws.Range("M24").Value = -11.111111111111
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -9.090909090909
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = 27.272727272727
$ws.Range("I25").Value = 19
$ws.Range("J25").Value = 20
$ws.Range("K25").Value = -5
$ws.Range("L25").Value = 11.764705882352
$ws.Range("M25").Value = 18.75
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -50
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 4
$ws.Range("K27").Value = 100
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0

# --- Numeric-count cell -> placeholder text "0" cell (style 14, shared string "0") ---
# Use Row 30 (same column) as a stable format-only source, never itself edited.
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0"
$ws.Range("C30").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C30").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("D30").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("E30").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C30").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("D30").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("E30").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("D30").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("E30").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("D30").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("E30").Copy()
$ws.Range("E29").PasteSpecial(-4122)

# --- Placeholder text cell -> numeric count cell (style 15) ---
# Use C36 (style 15) as a stable format-only source, never itself edited.
$ws.Range("D15").Value = 1
$ws.Range("C36").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("J15").Value = 1
$ws.Range("C36").Copy()
$ws.Range("J15").PasteSpecial(-4122)
$ws.Range("D26").Value = 3
$ws.Range("C36").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("J26").Value = 3
$ws.Range("C36").Copy()
$ws.Range("J26").PasteSpecial(-4122)

# --- Placeholder text cell -> numeric percent-change cell (style 16) ---
# Use K36 (style 16) as a stable format-only source, never itself edited.
$ws.Range("L14").Value = 0
$ws.Range("K36").Copy()
$ws.Range("L14").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("K36").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("K15").Value = 0
$ws.Range("K36").Copy()
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("M15").Value = 0
$ws.Range("K36").Copy()
$ws.Range("M15").PasteSpecial(-4122)
$ws.Range("L23").Value = 350
$ws.Range("K36").Copy()
$ws.Range("L23").PasteSpecial(-4122)
$ws.Range("E26").Value = -100
$ws.Range("K36").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("K26").Value = -66.666666666666
$ws.Range("K36").Copy()
$ws.Range("K26").PasteSpecial(-4122)

$excel.CutCopyMode = $false
